$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PartsList")

# Unify the "0603_cap"/"0603_CAP" footprint naming to "0603_CAP_SMALL"
$ws.Range("D2:D5").Value = "'0603_CAP_SMALL"

# Diode renames
$ws.Range("C9").Value = "'SCHOTTKY DIODE 0603"
$ws.Range("D9").Value = "'0603_DIODE-NSR20F30"

# USB connector footprint rename
$ws.Range("D10").Value = "'USB-MINI-B-HIR-UX60A-MB-5ST-SMALL"

# Resistor footprint renames
$ws.Range("D12:D18").Value = "'0603_res_SMALL"

# Column width adjustments (ColumnWidth is quantized internally to 1/6-character
# steps by this runtime, so inputs are pre-compensated to land as close as
# possible on the target serialized widths: 11.28515625, 23, 40, 22, 67.5703125)
$ws.Columns.Item(1).ColumnWidth = 10.5
$ws.Columns.Item(2).ColumnWidth = 22.166666666666668
$ws.Columns.Item(3).ColumnWidth = 39.166666666666664
$ws.Columns.Item(4).ColumnWidth = 21.166666666666668
$ws.Columns.Item(5).ColumnWidth = 66.66666666666667

# Remove explicit paper size (page setup)
$ws.PageSetup.PaperSize = $null
